# Daily attendance processing - 2026-01-08 10:37:50
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a value into a cell as LITERAL TEXT (never let Excel's
# autoparser turn a "NN.N%" looking string into a real percentage number),
# while preserving the destination cell's existing style/format (its "s"
# index does not change in the source diff). We stage the text in an
# otherwise-unused scratch cell (column J is blank throughout this sheet),
# force it to Text format there, copy it, and paste-special VALUES ONLY into
# the destination - a values-only paste carries the text over without
# dragging the scratch cell's own (Text) number format with it, so the
# target keeps its original style.
# ---------------------------------------------------------------------------
function Set-LiteralText($sheet, $cellRef, $text) {
    $scratch = $sheet.Range("J1")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $sheet.Range($cellRef).PasteSpecial(-4163)
    $scratch.Clear()
    $excel.CutCopyMode = 0
}

# ---------------------------------------------------------------------------
# 1) Summary box in column L, rows 6-10 (Recorded/Missing session counts and
#    the two coverage percentages recalculated after today's processing run)
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 219
$ws.Range("L7").Value = 27
Set-LiteralText $ws "L9" "68.9%"
Set-LiteralText $ws "L10" "75.9%"

# ---------------------------------------------------------------------------
# 2) "Recorded By" column (G) swaps order from
#    "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# ---------------------------------------------------------------------------
$swappedRows = @(8,9,10,12,14,15,17,18,34,35,36,38,40,41,43,44,60,61,62,64,66,67,69,70,86,87,88,90,92,93,95,96,112,113,114,116,118,119,121,122,138,139,140,142,144,145,147,148,164,167,170,174,191,194,197,201,218,221,224,228,245,248,251,255,272,275,278,282,299,302,305,309)
foreach ($r in $swappedRows) {
    $ws.Range("G$r").Value = "System, dnasr281@gmail.com"
}

# ---------------------------------------------------------------------------
# 3) Per-group breakdown table (columns K-S), rows 21-26 - attendance counts
#    recompute (one more Present, one fewer Absent) and the two % columns
#    refresh to match.
# ---------------------------------------------------------------------------
$ws.Range("O21").Value = 18
$ws.Range("P21").Value = 3
Set-LiteralText $ws "R21" "66.7%"
Set-LiteralText $ws "S21" "78.0%"

$ws.Range("O22").Value = 18
$ws.Range("P22").Value = 3
Set-LiteralText $ws "R22" "66.7%"
Set-LiteralText $ws "S22" "77.6%"

$ws.Range("O23").Value = 18
$ws.Range("P23").Value = 3
Set-LiteralText $ws "R23" "66.7%"
Set-LiteralText $ws "S23" "78.4%"

$ws.Range("O24").Value = 17
$ws.Range("P24").Value = 4
Set-LiteralText $ws "R24" "63.0%"
Set-LiteralText $ws "S24" "72.7%"

$ws.Range("O25").Value = 18
$ws.Range("P25").Value = 3
Set-LiteralText $ws "R25" "66.7%"
Set-LiteralText $ws "S25" "71.4%"

$ws.Range("O26").Value = 18
$ws.Range("P26").Value = 3
Set-LiteralText $ws "R26" "66.7%"
Set-LiteralText $ws "S26" "63.8%"

# ---------------------------------------------------------------------------
# 4) Six session rows that just got recorded (Not Recorded -> Recorded),
#    one per B1-group (B1D1, B1D2, B1E1, B1E2, B1F1, B1F2). Re-style the row
#    to the "Recorded" (green) look by pasting formats from an already
#    "Recorded" row, then fill in the recorder, the attendance fraction, and
#    the status text.
# ---------------------------------------------------------------------------
$recordedRows = @{
    178 = "17/23"
    205 = "26/30"
    232 = "14/26"
    259 = "21/28"
    286 = "21/26"
    313 = "18/29"
}

foreach ($r in $recordedRows.Keys) {
    $src = $ws.Range("A5:I5")
    $dst = $ws.Range("A" + $r + ":I" + $r)
    $src.Copy()
    $dst.PasteSpecial(-4122)
    $excel.CutCopyMode = 0

    $ws.Range("G$r").Value = "dnasr281@gmail.com"
    $ws.Range("H$r").Value = $recordedRows[$r]
    $ws.Range("I$r").Value = "Recorded"
}
